# Add two new inventory items (Hotdog Bun, Hamburger Bun) to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 10: Hotdog Bun
$ws.Cells.Item(10, 1).Value = 968149
$ws.Cells.Item(10, 2).Value = "Hotdog Bun"
$ws.Cells.Item(10, 3).Value = 1.23

# Row 11: Hamburger Bun
$ws.Cells.Item(11, 1).Value = 966705
$ws.Cells.Item(11, 2).Value = "Hamburger Bun"
$ws.Cells.Item(11, 3).Value = 1.23

# Match the style (left-aligned) used on column A for the existing rows
$ws.Range("A10:A11").HorizontalAlignment = -4131

# Update selection to reflect the last edited cell, as in the authored workbook
$ws.Range("C11").Select()
